$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Model description text (shared across all data rows)
$modelText = @"
MultiOutputRegressor(estimator=GridSearchCV(cv=5,
                                            estimator=Pipeline(steps=[('model',
                                                                       AdaBoostRegressor())]),
                                            param_grid={'model__learning_rate': [0.1,
                                                                                 0.5,
                                                                                 1.0],
                                                        'model__n_estimators': [50,
                                                                                100,
                                                                                150]},
                                            scoring='neg_mean_squared_error'))
"@

# New header cell F1 = "Modelo", styled the same as the other header cells (A1:E1)
$ws.Range("F1").Value = "Modelo"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Updated metric values (B,C,D) for each row, plus new model description (F)
$ws.Range("B2").Value = 0.5039322783415993
$ws.Range("C2").Value = 0.9899656381769188
$ws.Range("D2").Value = 0.5815779036458114
$ws.Range("F2").Value = $modelText

$ws.Range("B3").Value = 0.2393456043709881
$ws.Range("C3").Value = 0.9953219841434787
$ws.Range("D3").Value = 0.3827753431510709
$ws.Range("F3").Value = $modelText

$ws.Range("B4").Value = 0.3381587479270987
$ws.Range("C4").Value = 0.9934945776337086
$ws.Range("D4").Value = 0.4628283649783701
$ws.Range("F4").Value = $modelText

$ws.Range("B5").Value = 0.3851566315294757
$ws.Range("C5").Value = 0.9924052385526835
$ws.Range("D5").Value = 0.4846726063545287
$ws.Range("F5").Value = $modelText

$ws.Range("B6").Value = 0.6587938625413217
$ws.Range("C6").Value = 0.9806494985043637
$ws.Range("D6").Value = 0.6109652078713076
$ws.Range("F6").Value = $modelText

$ws.Range("B7").Value = 0.2479756927424193
$ws.Range("C7").Value = 0.9965729087833175
$ws.Range("D7").Value = 0.3768828815660853
$ws.Range("F7").Value = $modelText

$ws.Range("B8").Value = 0.1579536424331374
$ws.Range("C8").Value = 0.9983720054873471
$ws.Range("D8").Value = 0.3426586216972675
$ws.Range("F8").Value = $modelText

$ws.Range("B9").Value = 0.5213271728323167
$ws.Range("C9").Value = 0.9968946027779292
$ws.Range("D9").Value = 0.5974129538515721
$ws.Range("F9").Value = $modelText

$ws.Range("B10").Value = 0.1160246393334504
$ws.Range("C10").Value = 0.9978805046095817
$ws.Range("D10").Value = 0.2450065170048785
$ws.Range("F10").Value = $modelText

$ws.Range("B11").Value = 0.3228709807900418
$ws.Range("C11").Value = 0.9761451157290082
$ws.Range("D11").Value = 0.4408833004892758
$ws.Range("F11").Value = $modelText

$ws.Range("B12").Value = 0.04970959672068318
$ws.Range("C12").Value = 0.9985381255674414
$ws.Range("D12").Value = 0.1641097888004932
$ws.Range("F12").Value = $modelText

$ws.Range("B13").Value = 0.1019808664031013
$ws.Range("C13").Value = 0.999032665779645
$ws.Range("D13").Value = 0.223878916912976
$ws.Range("F13").Value = $modelText

$ws.Range("B14").Value = 0.09462228383192978
$ws.Range("C14").Value = 0.9987188726853985
$ws.Range("D14").Value = 0.2492067260330941
$ws.Range("F14").Value = $modelText

# Setting multi-line text via .Value auto-expands the row height; restore the
# default (non-custom) row height on the affected rows to match the original
# sheet formatting.
$ws.Range("A2:A14").EntireRow.AutoFit()
